# Final texts for 2.8 — deactivate the dragon reward in the daily reward
# feature (AB testing ARPDAU fix) on the "dailyLogin" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dailyLogin")
$ws.Activate()

# --- Row 9 (reward_6): drop its "pet" payout, take on the "hc" payout that
#     used to live on the now-removed reward_6b row. Priority flag -> 0.
$ws.Range("D9").Value = "hc"
$ws.Range("E9").Value = 3
$ws.Range("F9").ClearContents()
$ws.Range("G9").Value = 0

# --- Row 10 (reward_7): overwrite its sku label + day, take on the "pet"
#     payout that used to live on row 9. Priority flag -> 0 (dragon reward
#     deactivated).
$ws.Range("B10").Value = "reward_7"
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = "pet"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "pet_67"
$ws.Range("G10").Value = 0

# --- Row 11 (reward_14): overwrite its sku label + day, keep the "egg"
#     payout. Priority flag -> 0.
$ws.Range("B11").Value = "reward_14"
$ws.Range("C11").Value = 14
$ws.Range("D11").Value = "egg"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = "egg_dailyLogin"
$ws.Range("G11").Value = 0

# --- Rows 12:13 (reward_7b and the duplicate reward_14) are no longer
#     needed now their data moved up into rows 10/11 above — delete the
#     whole rows so everything below (the dragon modifiers table) shifts
#     up by two, matching the table/autofilter range shrinking from
#     A19:D39 to A17:D37.
$ws.Range("A12:A13").EntireRow.Delete()

# Cosmetic: match the author's last-saved selection.
$null = $ws.Range("J9").Select()
